# Add a third sheet ("Sheet1 Sample Invoice.pdf") that is a copy of the
# existing "Sheet1 Sample Invoice_32.pdf" sheet (same header row / data
# layout), but with its own Timestamp value in C3 - this gives the robot a
# dynamic selector/table input source to pick from.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Sheet1 Sample Invoice_32.pdf")

# Copy the source sheet and place the copy immediately after it (becomes
# the last sheet in the workbook).
$source.Copy([System.Reflection.Missing]::Value, $source)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Sheet1 Sample Invoice.pdf"

# Update the timestamp cell for this invoice sample to its own value.
$newSheet.Range("C3").Value = 44946.5077546296

$wb.Worksheets.Item(1).Select()
